# Update the "取得日時" (retrieved datetime) timestamps in the first sheet
# ("ランサーズ") from 2025-09-16 12:36:17 to 2025-09-16 12:45:13 for all
# data rows (rows 2 through 15, column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldValue = "2025-09-16 12:36:17"
$newValue = "2025-09-16 12:45:13"

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
